# "code changes are done for watchlist and notification"
#
# Target sheet is "Test Cases" (the active / tab-selected sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Runmode column (C) for the watchlist test cases (rows 2-8): "Y" -> "N"
# (they are temporarily turned off while the watchlist/notification code
# changes are verified).
$ws.Range("C2:C8").Value = "N"

# Results column (D) for rows 9-11: "SKIP" -> "FAIL"
$ws.Range("D9:D11").Value = "FAIL"

# Column A is widened / best-fit to accommodate longer TCID text.
$ws.Columns.Item(1).ColumnWidth = 39.6

# Leave the selection on C19, matching where the author finished editing.
[void]$ws.Range("C19").Select()
